$d = $word.ActiveDocument
$t = $d.Tables.Item(5)

# Step 1: seed each target (currently empty) cell with a unique placeholder token
$t.Cell(2,2).Range.Text = "@@CELL00@@"
$t.Cell(2,3).Range.Text = "@@CELL01@@"
$t.Cell(2,4).Range.Text = "@@CELL02@@"
$t.Cell(3,2).Range.Text = "@@CELL03@@"
$t.Cell(3,3).Range.Text = "@@CELL04@@"
$t.Cell(3,4).Range.Text = "@@CELL05@@"
$t.Cell(4,2).Range.Text = "@@CELL06@@"
$t.Cell(4,3).Range.Text = "@@CELL07@@"
$t.Cell(4,4).Range.Text = "@@CELL08@@"
$t.Cell(5,2).Range.Text = "@@CELL09@@"
$t.Cell(5,3).Range.Text = "@@CELL10@@"
$t.Cell(5,4).Range.Text = "@@CELL11@@"
$t.Cell(6,2).Range.Text = "@@CELL12@@"
$t.Cell(6,3).Range.Text = "@@CELL13@@"
$t.Cell(6,4).Range.Text = "@@CELL14@@"
$t.Cell(7,2).Range.Text = "@@CELL15@@"
$t.Cell(7,3).Range.Text = "@@CELL16@@"
$t.Cell(7,4).Range.Text = "@@CELL17@@"
$t.Cell(8,2).Range.Text = "@@CELL18@@"
$t.Cell(8,3).Range.Text = "@@CELL19@@"
$t.Cell(8,4).Range.Text = "@@CELL20@@"
$t.Cell(9,2).Range.Text = "@@CELL21@@"
$t.Cell(9,3).Range.Text = "@@CELL22@@"
$t.Cell(9,4).Range.Text = "@@CELL23@@"
$t.Cell(10,2).Range.Text = "@@CELL24@@"
$t.Cell(10,3).Range.Text = "@@CELL25@@"
$t.Cell(10,4).Range.Text = "@@CELL26@@"
$t.Cell(11,2).Range.Text = "@@CELL27@@"
$t.Cell(11,3).Range.Text = "@@CELL28@@"
$t.Cell(11,4).Range.Text = "@@CELL29@@"
$t.Cell(12,2).Range.Text = "@@CELL30@@"
$t.Cell(12,3).Range.Text = "@@CELL31@@"
$t.Cell(12,4).Range.Text = "@@CELL32@@"
$t.Cell(13,2).Range.Text = "@@CELL33@@"
$t.Cell(13,3).Range.Text = "@@CELL34@@"
$t.Cell(13,4).Range.Text = "@@CELL35@@"
$t.Cell(14,2).Range.Text = "@@CELL36@@"
$t.Cell(14,3).Range.Text = "@@CELL37@@"
$t.Cell(14,4).Range.Text = "@@CELL38@@"
$t.Cell(15,2).Range.Text = "@@CELL39@@"
$t.Cell(15,3).Range.Text = "@@CELL40@@"
$t.Cell(15,4).Range.Text = "@@CELL41@@"
$t.Cell(16,2).Range.Text = "@@CELL42@@"
$t.Cell(16,3).Range.Text = "@@CELL43@@"
$t.Cell(16,4).Range.Text = "@@CELL44@@"
$t.Cell(17,2).Range.Text = "@@CELL45@@"
$t.Cell(17,3).Range.Text = "@@CELL46@@"
$t.Cell(17,4).Range.Text = "@@CELL47@@"
$t.Cell(18,2).Range.Text = "@@CELL48@@"
$t.Cell(18,3).Range.Text = "@@CELL49@@"
$t.Cell(18,4).Range.Text = "@@CELL50@@"
$t.Cell(19,2).Range.Text = "@@CELL51@@"
$t.Cell(19,3).Range.Text = "@@CELL52@@"
$t.Cell(19,4).Range.Text = "@@CELL53@@"
$t.Cell(20,2).Range.Text = "@@CELL54@@"
$t.Cell(20,3).Range.Text = "@@CELL55@@"
$t.Cell(20,4).Range.Text = "@@CELL56@@"
$t.Cell(21,2).Range.Text = "@@CELL57@@"
$t.Cell(21,3).Range.Text = "@@CELL58@@"
$t.Cell(21,4).Range.Text = "@@CELL59@@"

# Step 2: Find/Replace each placeholder with the real value, applying the
# run formatting (Times New Roman / bold / 14pt) used throughout this table
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Times New Roman"
$find.Replacement.Font.NameBi = "Times New Roman"
$find.Replacement.Font.Bold = $true
$find.Replacement.Font.Size = 14
$find.Execute("@@CELL00@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.6.0", 2)
$find.Execute("@@CELL01@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.224", 2)
$find.Execute("@@CELL02@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.246", 2)
$find.Execute("@@CELL03@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255", 2)
$find.Execute("@@CELL04@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL05@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.246", 2)
$find.Execute("@@CELL06@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.0.0", 2)
$find.Execute("@@CELL07@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.252.0", 2)
$find.Execute("@@CELL08@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.246", 2)
$find.Execute("@@CELL09@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.4.0", 2)
$find.Execute("@@CELL10@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL11@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.249", 2)
$find.Execute("@@CELL12@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.252", 2)
$find.Execute("@@CELL13@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL14@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.249", 2)
$find.Execute("@@CELL15@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.5.0", 2)
$find.Execute("@@CELL16@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL17@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.245", 2)
$find.Execute("@@CELL18@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.248", 2)
$find.Execute("@@CELL19@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL20@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.245", 2)
$find.Execute("@@CELL21@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.4.0", 2)
$find.Execute("@@CELL22@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL23@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.245", 2)
$find.Execute("@@CELL24@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.252", 2)
$find.Execute("@@CELL25@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL26@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.242", 2)
$find.Execute("@@CELL27@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.0.0", 2)
$find.Execute("@@CELL28@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL29@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.242", 2)
$find.Execute("@@CELL30@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.5.0", 2)
$find.Execute("@@CELL31@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL32@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.254", 2)
$find.Execute("@@CELL33@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.248", 2)
$find.Execute("@@CELL34@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL35@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.254", 2)
$find.Execute("@@CELL36@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.4.0", 2)
$find.Execute("@@CELL37@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL38@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.254", 2)
$find.Execute("@@CELL39@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.244", 2)
$find.Execute("@@CELL40@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL41@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.241", 2)
$find.Execute("@@CELL42@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.6.0", 2)
$find.Execute("@@CELL43@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.224", 2)
$find.Execute("@@CELL44@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.241", 2)
$find.Execute("@@CELL45@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.5.0", 2)
$find.Execute("@@CELL46@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.0", 2)
$find.Execute("@@CELL47@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.250", 2)
$find.Execute("@@CELL48@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.244", 2)
$find.Execute("@@CELL49@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL50@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.250", 2)
$find.Execute("@@CELL51@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.240", 2)
$find.Execute("@@CELL52@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.252", 2)
$find.Execute("@@CELL53@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.253", 2)
$find.Execute("@@CELL54@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.6.0", 2)
$find.Execute("@@CELL55@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.255.224", 2)
$find.Execute("@@CELL56@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.253", 2)
$find.Execute("@@CELL57@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.0.0", 2)
$find.Execute("@@CELL58@@", $false, $false, $false, $false, $false, $true, 1, $false, "255.255.252.0", 2)
$find.Execute("@@CELL59@@", $false, $false, $false, $false, $false, $true, 1, $false, "16.128.255.253", 2)
